# Generate Report for Handback
# Removes the stale "eb5361da-e598-49e1-a781-2298c27002b2" report row from every
# sheet (Overview / zh-cn / de-de) and refreshes the handoff/handback
# timestamps for the remaining "bfc56466-f424-4c71-a9e5-4b645e843490" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop every hyperlink on the sheet (only way this host exposes bulk removal)
# so we can rebuild just the ones that should survive the row delete.
$wsOverview.Range("A1").Hyperlinks.Delete()

# Remove the whole eb5361da... row (row 3); row 2 (bfc56466...) stays put.
$wsOverview.Rows.Item(3).Delete()

# Re-create the hyperlink that belongs to the surviving row.
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2aca258f9e1e0f31e6c090b2896561b71794e03/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A1").Hyperlinks.Delete()

$wsZh.Rows.Item(3).Delete()

# Refresh the handoff / handback timestamps for the remaining row.
$wsZh.Range("E2").Value2 = "2016-03-24 00:52:47"
$wsZh.Range("H2").Value2 = "2016-03-24 00:53:10"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2aca258f9e1e0f31e6c090b2896561b71794e03/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d334afd6da1c26c86212a43bf72de19f6bb40da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4890fdab488a4a03b7973da9ef7a52a905d43c3c/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/80eb455c6054d4ac1faaa47d25ef5460f964e8b2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A1").Hyperlinks.Delete()

$wsDe.Rows.Item(3).Delete()

$wsDe.Range("E2").Value2 = "2016-03-24 00:52:51"
$wsDe.Range("H2").Value2 = "2016-03-24 00:53:16"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2aca258f9e1e0f31e6c090b2896561b71794e03/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b6f8c878dd73f33eff0fc18d2062985d2eca290/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ee88f1c2c5370f97814320690e9167366f5add9d/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/db4bfda931f04058a061f49c8c6274d0fb8dc047/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf"
) | Out-Null

Write-Output "done"
